$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Sent By" column header
$ws.Range("G1").Value = "Sent By"

# Row 2 (ADM001) updates
$ws.Range("D2").Value = "04 February 2025, 12:00 AM"
$ws.Range("E2").Value = "Hello, your attendance has been recorded."
$ws.Range("F2").Value = "Pending"
$ws.Range("G2").Value = "John Smith"

# Row 3 (ADM002) updates
$ws.Range("D3").Value = "04 February 2025, 12:00 AM"
$ws.Range("E3").Value = "Reminder: Your exam is scheduled for tomorrow."
$ws.Range("F3").Value = "Pending"
$ws.Range("G3").Value = "John Smith"
